# ABBYY_IssueList.xlsx update:
#  - rename the existing sheet to "Original"
#  - duplicate it as a new sheet named "Extra Crispy"
#  - on "Extra Crispy", replace the Script text (column C) for every data row
#    with the updated search-rectangle variants, and flip E9 (TwelveFourPF)
#    from FALSE to TRUE
#  - re-select sheet/ranges to match the authored workbook

$wb = $excel.ActiveWorkbook

# --- Step 1: rename the original (only) sheet ---------------------------
$original = $wb.Worksheets.Item(1)
$original.Name = "Original"

# --- Step 2: duplicate it to create "Extra Crispy" -----------------------
$original.Copy([System.Reflection.Missing]::Value, $original)
$extra = $wb.Worksheets.Item(2)
$extra.Name = "Extra Crispy"

# --- Step 3: update the Script column text on "Extra Crispy" -------------
$extra.Range("C2").Value = 'Imagefound(text:"Episode", SearchRectangle: [231,209,1896,247], waitFor: 1)'
$extra.Range("C3").Value = 'Imagefound(TEXT:"Application Access Menu", validWords:"*")'
$extra.Range("C4").Value = 'Imagefound(TEXT:"ZZZeggplant, IPInfection",ValidCharacters:"A".."Z" &&& "a".."z" &&& "," &&& ".",IgnoreSpaces:yes,waitfor:0,textDifference:1,ignorenewlines:yes ,searchRectangle:[465,194,573,1078])'
$extra.Range("C5").Value = 'Imagefound(text:"In", waitFor:0,caseSensitive:"yes", SearchRectangle: [655,347,1265,652])'
$extra.Range("C6").Value = 'Imagefound(TEXT:"ZZZeggplant, IPwardnurse",ValidCharacters:"A".."Z" &&& "a".."z" &&& "," &&& ".",IgnoreSpaces:yes,waitfor:0,textDifference:1,ignorenewlines:yes ,searchRectangle:[465,194,573,1078])'
$extra.Range("C7").Value = 'imageFound(text: "patient", searchRectangle:[1240,143,1915,1034])'
$extra.Range("C8").Value = 'imageFound(dpi:"144", SearchRectangle:[230,240,1919,279], text:"Hide Add''l Visits", TextDifference:"2")'
$extra.Range("C9").Value = 'Imagefound(SearchRectangle:[435,155,1482,972], text:"Arrange Views")'

# TwelveFourPF for row 9 goes from FALSE to TRUE
$extra.Range("E9").Value = $true

# --- Step 4: selection bookkeeping ---------------------------------------
# "Extra Crispy" ends up with E2:E9 highlighted (active cell at the bottom)
$extra.Range("E2:E9").Select()

# Put the focus back on "Original" (first / visible tab) and reset its
# selection back to the top-left corner.
$original.Activate()
$original.Range("A1").Select()
